$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.222.80"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "2.005.53"
$ws.Range("E3").Value = "  +2.12%  "

$ws.Range("E4").Value = "  +0.04%  "

$cell = $ws.Range("D5")
$cell.Value = "'246.48"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("E6").Value = "  +1.32%  "

$cell = $ws.Range("D7")
$cell.Value = "'60.00"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -1.15%  "

$ws.Range("E8").Value = "  -0.04%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.386"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.16%  "

$cell = $ws.Range("D10")
$cell.Value = "'0.0807"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("E11").Value = "  +0.72%  "

$cell = $ws.Range("D12")
$cell.Value = "'15.13"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +6.15%  "

$cell = $ws.Range("D13")
$cell.Value = "'22.50"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").Value = "2.295.94"
$ws.Range("E14").Value = "  +2.04%  "

$cell = $ws.Range("D15")
$cell.Value = "'0.844"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "

$cell = $ws.Range("D16")
$cell.Value = "'5.44"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.49%  "

$ws.Range("D17").Value = "2.003.14"
$ws.Range("E17").Value = "  +2.33%  "

$ws.Range("D18").Value = "37.139.48"
$ws.Range("E18").Value = "  +1.55%  "

$cell = $ws.Range("D19")
$cell.Value = "'70.25"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("E20").Value = "  +1.41%  "

$cell = $ws.Range("D21")
$cell.Value = "'5.20"
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.Value = "'230.17"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  -0.08%  "

$cell = $ws.Range("D25")
$cell.Value = "'2.36"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("E26").Value = "  +2.16%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D27")
$cell.Value = "'164.49"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D28")
$cell.Value = "'0.139"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -5.07%  "

$cell = $ws.Range("D29")
$cell.Value = "'19.65"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

$cell = $ws.Range("D30")
$cell.Value = "'1.40"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +15.78%  "

$ws.Range("E31").Value = "  +1.00%  "

$cell = $ws.Range("D32")
$cell.Value = "'4.81"
$cell.Style = "Normal"

$ws.Range("E33").Value = "  +5.78%  "

$ws.Range("E34").Value = "  +1.05%  "

$cell = $ws.Range("D35")
$cell.Value = "'2.41"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +5.02%  "

$ws.Range("E36").Value = "  +0.14%  "

$cell = $ws.Range("D37")
$cell.Value = "'1.81"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.40%  "

$cell = $ws.Range("D38")
$cell.Value = "'3.34"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -6.64%  "

$cell = $ws.Range("D39")
$cell.Value = "'5.37"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -4.29%  "

$cell = $ws.Range("D40")
$cell.Value = "'0.0987"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  +1.48%  "

$cell = $ws.Range("D44")
$cell.Value = "'16.63"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +2.97%  "

$cell = $ws.Range("D45")
$cell.Value = "'91.50"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +3.33%  "

$ws.Range("D46").Value = "1.370.77"
$ws.Range("E46").Value = "  -0.05%  "

$cell = $ws.Range("D47")
$cell.Value = "'1.05"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "

$cell = $ws.Range("D48")
$cell.Value = "'7.35"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.69%  "

$cell = $ws.Range("D49")
$cell.Value = "'2.09"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +13.54%  "

$cell = $ws.Range("D50")
$cell.Value = "'47.09"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +5.71%  "

$ws.Range("E51").Value = "  -0.39%  "
